$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for the wxGroup column
$ws.Range("D1").Value = "wxGroup"

# Fill the wxGroup value ("days") for all data rows (2 through 101)
$ws.Range("D2:D101").Value = "days"

# Update the selection (also clears any frozen/top-left scroll position)
$ws.Range("G13").Select()
